# Apply scraper update to top_remaining_KY.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. LAST SCRAPE DATE moved forward for every data row (F2:F59).
#    Force the cells to remain plain text (the source value looks like a
#    date, so Excel would otherwise silently convert it to a date serial).
$ws.Range("F2:F59").NumberFormat = "@"
$ws.Range("F2:F59").Value = '2019-03-12'
$ws.Range("F2:F59").ClearFormats()

# 2. Rows 14/15 (10X THE CASH group) - game order / prize counts refreshed
$ws.Range("C14").Value = 'Holiday Gold $10,000'
$ws.Range("D14").Value = 698
$ws.Range("E14").Value = 2

$ws.Range("C15").Value = 'Find $100'
$ws.Range("D15").Value = 707
$ws.Range("E15").Value = 5

# 3. Row 19 (Crossword) - prize count refreshed
$ws.Range("E19").Value = 10

# 4. Rows 31/32 ($5.00 Games group) - game order / prize counts refreshed
$ws.Range("C31").Value = 'Giant Jumbo Bucks'
$ws.Range("D31").Value = 674
$ws.Range("E31").Value = 3

$ws.Range("C32").Value = 'Deluxe Crossword Doubler'
$ws.Range("D32").Value = 533
$ws.Range("E32").Value = 2

# 5. Rows 40/41 ($10.00 Games group) - game order / prize counts refreshed
$ws.Range("C40").Value = '$10,000,000 Cash Explosion'
$ws.Range("D40").Value = 530
$ws.Range("E40").Value = 1

$ws.Range("C41").Value = 'Holiday Gold $150,000'
$ws.Range("D41").Value = 700
$ws.Range("E41").Value = 2

# 6. Row 46 ($25, $50, or $100) - prize count refreshed
$ws.Range("E46").Value = 2251
